$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'31.267.48"
$ws.Range("E2").Value = '  +1.74%  '

# Row 3
$ws.Range("D3").Value = "'1.956.18"
$ws.Range("E3").Value = '  +0.42%  '

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").Value = "'246.44"
$ws.Range("E5").Value = '  -0.72%  '

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.09%  '

# Row 7
$ws.Range("D7").Value = "'0.4883"
$ws.Range("E7").Value = '  +1.38%  '

# Row 8
$ws.Range("D8").Value = "'0.2977"
$ws.Range("E8").Value = '  +1.17%  '

# Row 9
$ws.Range("D9").Value = "'0.06848"
$ws.Range("E9").Value = '  +0.32%  '

# Row 10
$ws.Range("D10").Value = "'19.44"
$ws.Range("E10").Value = '  -0.10%  '

# Row 11
$ws.Range("D11").Value = "'107.36"
$ws.Range("E11").Value = '  -5.11%  '

# Row 12
$ws.Range("D12").Value = "'0.07738"
$ws.Range("E12").Value = '  +1.39%  '

# Row 13
$ws.Range("D13").Value = "'1.923.06"
$ws.Range("E13").Value = '  -1.33%  '

# Row 14
$ws.Range("D14").Value = "'5.458"
$ws.Range("E14").Value = '  -2.27%  '

# Row 15
$ws.Range("D15").Value = "'0.7111"
$ws.Range("E15").Value = '  +3.36%  '

# Row 16
$ws.Range("D16").Value = "'283.61"
$ws.Range("E16").Value = '  -5.41%  '

# Row 17
$ws.Range("D17").Value = "'31.157.68"
$ws.Range("E17").Value = '  +1.40%  '

# Row 18
$ws.Range("D18").Value = "'0.000007759"
$ws.Range("E18").Value = '  +0.64%  '

# Row 19
$ws.Range("D19").Value = "'13.20"
$ws.Range("E19").Value = '  -0.45%  '

# Row 20
$ws.Range("D20").Value = "'1.0000"
$ws.Range("E20").Value = '  +0.06%  '

# Row 21
$ws.Range("D21").Value = "'2.190.00"
$ws.Range("E21").Value = '  -0.48%  '

# Row 22
$ws.Range("D22").Value = "'5.502"
$ws.Range("E22").Value = '  -2.88%  '

# Row 23
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = '  +0.19%  '

# Row 24
$ws.Range("D24").Value = "'6.569"
$ws.Range("E24").Value = '  -0.49%  '

# Row 25
$ws.Range("D25").Value = "'9.869"
$ws.Range("E25").Value = '  +1.75%  '

# Row 26
$ws.Range("D26").Value = "'169.61"
$ws.Range("E26").Value = '  +1.03%  '

# Row 27
$ws.Range("D27").Value = "'20.31"
$ws.Range("E27").Value = '  -1.49%  '

# Row 28
$ws.Range("D28").Value = "'2.196"
$ws.Range("E28").Value = '  +1.77%  '

# Row 29
$ws.Range("D29").Value = "'0.1053"
$ws.Range("E29").Value = '  -2.53%  '

# Row 30
$ws.Range("D30").Value = "'1.435"
$ws.Range("E30").Value = '  +0.43%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'1.595"
$ws.Range("E31").Value = '  +0.16%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'4.652"
$ws.Range("E32").Value = '  -0.10%  '

# Row 33
$ws.Range("D33").Value = "'4.445"
$ws.Range("E33").Value = '  +1.22%  '

# Row 34
$ws.Range("D34").Value = "'0.04984"
$ws.Range("E34").Value = '  -1.69%  '

# Row 35
$ws.Range("D35").Value = "'0.7600"
$ws.Range("E35").Value = '  -1.64%  '

# Row 36
$ws.Range("D36").Value = "'1.166"
$ws.Range("E36").Value = '  -0.10%  '

# Row 37
$ws.Range("D37").Value = "'2.733"
$ws.Range("E37").Value = '  +0.01%  '

# Row 38
$ws.Range("D38").Value = "'0.02039"
$ws.Range("E38").Value = '  -2.06%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = "'2.220"
$ws.Range("E39").Value = '  +8.75%  '

# Row 40
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = "'2.705"
$ws.Range("E40").Value = '  -0.27%  '

# Row 41
$ws.Range("D41").Value = "'6.420"
$ws.Range("E41").Value = '  +8.57%  '

# Row 42
$ws.Range("D42").Value = "'0.4534"
$ws.Range("E42").Value = '  +0.97%  '

# Row 43
$ws.Range("D43").Value = "'109.50"
$ws.Range("E43").Value = '  -1.99%  '

# Row 44
$ws.Range("D44").Value = "'0.8824"
$ws.Range("E44").Value = '  +0.61%  '

# Row 45
$ws.Range("D45").Value = "'71.86"
$ws.Range("E45").Value = '  +1.03%  '

# Row 46
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = '  -0.16%  '

# Row 47
$ws.Range("D47").Value = "'7.803"
$ws.Range("E47").Value = '  +5.39%  '

# Row 48
$ws.Range("D48").Value = "'0.2630"
$ws.Range("E48").Value = '  +2.06%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'9.356"
$ws.Range("E49").Value = '  -1.71%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = "'0.1262"
$ws.Range("E50").Value = '  +1.00%  '

# Row 51
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = "'955.73"
$ws.Range("E51").Value = '  +5.63%  '
